# chore: update Sheets via scheduled runner
# Refreshes cached market-board price/profit figures (columns H-N) for a
# handful of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1666675.4
$ws.Range("I6").Value = 1666675.4
$ws.Range("K6").Value = 5000026.199999999
$ws.Range("M6").Value = -4999914.199999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 72.25
$ws.Range("I9").Value = 83
$ws.Range("J9").Value = 40
$ws.Range("K9").Value = 83
$ws.Range("L9").Value = 40
$ws.Range("M9").Value = 86
$ws.Range("N9").Value = -378

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 123.71429
$ws.Range("I33").Value = 123.71429
$ws.Range("K33").Value = 123.71429
$ws.Range("M33").Value = 105.28571

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1518.75
$ws.Range("I38").Value = 50
$ws.Range("J38").Value = 3966.6667
$ws.Range("K38").Value = 150
$ws.Range("L38").Value = 11900.0001
$ws.Range("M38").Value = 222
$ws.Range("N38").Value = -12644.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5000
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -4516
$ws.Range("N51").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3874.5715
$ws.Range("I2").Value = 2374
$ws.Range("K2").Value = 2374
$ws.Range("M2").Value = -2261

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 18000
$ws.Range("J24").Value = 18000
$ws.Range("L24").Value = 18000
$ws.Range("N24").Value = -18748

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2868
$ws.Range("I32").Value = 2846
$ws.Range("K32").Value = 2846
$ws.Range("M32").Value = -2559

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3612.3333
$ws.Range("I63").Value = 4474.5
$ws.Range("K63").Value = 4474.5
$ws.Range("M63").Value = -3788.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3612.3333
$ws.Range("I66").Value = 4474.5
$ws.Range("K66").Value = 22372.5
$ws.Range("M66").Value = -18940.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H100").Value = 18000
$ws.Range("J100").Value = 18000
$ws.Range("L100").Value = 18000
$ws.Range("N100").Value = -20164

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3874.5715
$ws.Range("I116").Value = 2374
$ws.Range("K116").Value = 2374
$ws.Range("M116").Value = -80

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 58285.715
$ws.Range("J124").Value = 58285.715
$ws.Range("L124").Value = 58285.715
$ws.Range("N124").Value = -68105.715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 49000
$ws.Range("J125").Value = 49000
$ws.Range("L125").Value = 49000
$ws.Range("N125").Value = -58840

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 40214.5
$ws.Range("J134").Value = 40214.5
$ws.Range("L134").Value = 40214.5
$ws.Range("N134").Value = -50354.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3874.5715
$ws.Range("I3").Value = 2374
$ws.Range("K3").Value = 2374
$ws.Range("M3").Value = -2260

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 398.4
$ws.Range("I11").Value = 248.25
$ws.Range("K11").Value = 248.25
$ws.Range("M11").Value = -108.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1089.6666
$ws.Range("I20").Value = 1280
$ws.Range("J20").Value = 899.3333
$ws.Range("K20").Value = 1280
$ws.Range("L20").Value = 899.3333
$ws.Range("M20").Value = -1033
$ws.Range("N20").Value = -1393.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 43.666668
$ws.Range("I37").Value = 43.666668
$ws.Range("K37").Value = 43.666668
$ws.Range("M37").Value = 93.333332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3000
$ws.Range("I86").Value = 3000
$ws.Range("K86").Value = 3000
$ws.Range("M86").Value = -1877

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3000
$ws.Range("I89").Value = 3000
$ws.Range("K89").Value = 15000
$ws.Range("M89").Value = -9384

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2735.75
$ws.Range("I105").Value = 2735.75
$ws.Range("K105").Value = 2735.75
$ws.Range("M105").Value = -988.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2170.6667
$ws.Range("I6").Value = 2341.3333
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 2341.3333
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = -2228.3333
$ws.Range("N6").Value = -2226

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6293.625
$ws.Range("I31").Value = 1851.7142
$ws.Range("K31").Value = 1851.7142
$ws.Range("M31").Value = -1556.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6293.625
$ws.Range("I34").Value = 1851.7142
$ws.Range("K34").Value = 1851.7142
$ws.Range("M34").Value = -1649.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 8392.333000000001
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 8392.333000000001
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 52081
$ws.Range("J93").Value = 43333
$ws.Range("L93").Value = 43333
$ws.Range("N93").Value = -47077

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1515.1818
$ws.Range("I122").Value = 1140.7778
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 3422.3334
$ws.Range("L122").Value = 9600
$ws.Range("M122").Value = -972.3334000000004
$ws.Range("N122").Value = -14500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 12.933333
$ws.Range("I2").Value = 14.083333
$ws.Range("K2").Value = 84.49999800000001
$ws.Range("M2").Value = 28.50000199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 70.57143000000001
$ws.Range("I7").Value = 68.8
$ws.Range("J7").Value = 75
$ws.Range("K7").Value = 206.4
$ws.Range("L7").Value = 225
$ws.Range("M7").Value = -94.39999999999998
$ws.Range("N7").Value = -449

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 150
$ws.Range("I34").Value = 150
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 450
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -366
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1400
$ws.Range("I68").Value = 800
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 2400
$ws.Range("L68").Value = 6000
$ws.Range("M68").Value = -1589
$ws.Range("N68").Value = -7622

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1400
$ws.Range("I71").Value = 800
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 7200
$ws.Range("L71").Value = 18000
$ws.Range("M71").Value = -3144
$ws.Range("N71").Value = -26112

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1957.2858
$ws.Range("J129").Value = 1950.3334
$ws.Range("L129").Value = 5851.0002
$ws.Range("N129").Value = -15851.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1861.7059
$ws.Range("J131").Value = 2299.5833
$ws.Range("L131").Value = 6898.749899999999
$ws.Range("N131").Value = -16978.7499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10701.6
$ws.Range("I70").Value = 6000
$ws.Range("J70").Value = 11877
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 11877
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -12417

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 10701.6
$ws.Range("I73").Value = 6000
$ws.Range("J73").Value = 11877
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 11877
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -13749

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 1000
$ws.Range("K97").Value = 1000
$ws.Range("M97").Value = -504

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
